# Update "Pais" COVID dashboard sheet: refresh data for a handful of
# countries, re-sort the table by "Casos totales" (column B) descending
# (ties keep their previous relative order), and bump the "last updated"
# timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Datos actualizados..." timestamp in A1.
$ws.Range("A1").Value2 = "Datos actualizados a 13 de Abril de 2020 a las 06:22"

# 2) Push fresh totals for the countries that changed in this refresh.
#    (row numbers are the *current*, pre-sort positions in the sheet)
#    Kazajistan: row 75
$ws.Cells.Item(75, 2).Value2 = 961
$ws.Cells.Item(75, 3).Value2 = 10
$ws.Cells.Item(75, 5).Value2 = 852

#    Uruguay: row 96
$ws.Cells.Item(96, 5).Value2 = 241
$ws.Cells.Item(96, 6).Value2 = 15
$ws.Cells.Item(96, 7).Value2 = 1
$ws.Cells.Item(96, 8).Value2 = 8

#    Honduras: row 98
$ws.Cells.Item(98, 2).Value2 = 397
$ws.Cells.Item(98, 3).Value2 = 4
$ws.Cells.Item(98, 5).Value2 = 365

#    Paraguay: row 126 (before re-sort)
$ws.Cells.Item(126, 2).Value2 = 147
$ws.Cells.Item(126, 3).Value2 = 13
$ws.Cells.Item(126, 5).Value2 = 119
$ws.Cells.Item(126, 6).Value2 = 1
$ws.Cells.Item(126, 8).Value2 = 6

#    Mongolia: row 178 (before re-sort)
$ws.Cells.Item(178, 2).Value2 = 17
$ws.Cells.Item(178, 3).Value2 = 1
$ws.Cells.Item(178, 5).Value2 = 13

# 3) The countries table (A4:H216) is kept sorted by "Casos totales"
#    (column B) descending. Paraguay's and Mongolia's new totals move
#    them ahead of some neighbours, so re-apply that ordering to the two
#    affected neighbourhoods by swapping the now out-of-order rows.

# Paraguay (row 126, B=147) now outranks Guadalupe (row 124, B=143) and
# Brunei (row 125, B=136): rotate rows 124-126 so Paraguay lands on top
# while Guadalupe and Brunei keep their relative order and data.
$cols = 1..8
$row124 = @{}
$row125 = @{}
$row126 = @{}
foreach ($c in $cols) {
    $row124[$c] = $ws.Cells.Item(124, $c).Value2
    $row125[$c] = $ws.Cells.Item(125, $c).Value2
    $row126[$c] = $ws.Cells.Item(126, $c).Value2
}
foreach ($c in $cols) {
    $ws.Cells.Item(124, $c).Value2 = $row126[$c]
    $ws.Cells.Item(125, $c).Value2 = $row124[$c]
    $ws.Cells.Item(126, $c).Value2 = $row125[$c]
}

# Mongolia (row 178, B=17) now outranks Fiyi (row 176, B=16) and Namibia
# (row 177, B=16): rotate rows 176-178 the same way.
$row176 = @{}
$row177 = @{}
$row178 = @{}
foreach ($c in $cols) {
    $row176[$c] = $ws.Cells.Item(176, $c).Value2
    $row177[$c] = $ws.Cells.Item(177, $c).Value2
    $row178[$c] = $ws.Cells.Item(178, $c).Value2
}
foreach ($c in $cols) {
    $ws.Cells.Item(176, $c).Value2 = $row178[$c]
    $ws.Cells.Item(177, $c).Value2 = $row176[$c]
    $ws.Cells.Item(178, $c).Value2 = $row177[$c]
}
